$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit performs a cyclic rotation of three observation records
# (rows 3, 4 and 5): the data that used to be in row 5 moves to row 3,
# the data that used to be in row 3 moves to row 4, and the data that
# used to be in row 4 moves to row 5. Only a subset of columns (A, B, D,
# E, F, G, H, Q, R, Y, AA) actually carries information that differs
# between these rows; everything else (C, I, P, S, T, U, V, W, Z, AB,
# AD, AE, AG, AT, AW, AX, AY) is identical across the three rows and is
# therefore left untouched.

function Set-RowData {
    param($ws, $row, $data)

    $ws.Range("A$row").Value = $data.A
    $ws.Range("B$row").Value = $data.B
    $ws.Range("D$row").Value = $data.D
    $ws.Range("E$row").Value = $data.E
    $ws.Range("F$row").Value = $data.F
    $ws.Range("G$row").Value = $data.G
    $ws.Range("H$row").Value = $data.H
    $ws.Range("Q$row").Value = $data.Q
    $ws.Range("R$row").Value = $data.R

    # Y/AA hold plain date-like text (e.g. "2022-04-24"); Excel's Value
    # setter would otherwise auto-convert such strings into date serial
    # numbers. Force the cell to stay text, write the value, then reset
    # the style so no spurious numFmt/style index is left behind.
    $ws.Range("Y$row").NumberFormat = "@"
    $ws.Range("Y$row").Value = $data.Y
    $ws.Range("Y$row").Style = "Normal"

    $ws.Range("AA$row").NumberFormat = "@"
    $ws.Range("AA$row").Value = $data.AA
    $ws.Range("AA$row").Style = "Normal"
}

$newRow3 = @{
    A = 105030033; B = 78098; D = "NT"; E = 6453
    F = "Vedskivlav"; G = "Hertelidea botryosa"; H = "(Fr.) Printzen & Kantvilas"
    Q = 503637.7373408998; R = 6838237.474158124
    Y = "2022-04-24"; AA = "2022-04-24"
}

$newRow4 = @{
    A = 105030138; B = 77177; D = "NT"; E = 353
    F = "Dvärgbägarlav"; G = "Cladonia parasitica"; H = "(Hoffm.) Hoffm."
    Q = 503511.5508635575; R = 6838417.747975093
    Y = "2022-04-22"; AA = "2022-04-22"
}

$newRow5 = @{
    A = 105030139; B = 90653; D = "LC"; E = 4364
    F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."
    Q = 503513.9504220288; R = 6838395.91245344
    Y = "2022-04-22"; AA = "2022-04-22"
}

Set-RowData $ws 3 $newRow3
Set-RowData $ws 4 $newRow4
Set-RowData $ws 5 $newRow5
